$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.450.82"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.62"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.03"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "3.975.78"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.24"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "3.400.88"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "61.475.04"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.70"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.97"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.77"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.37"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("E26").Value = "  +8.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "167.71"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "3.427.19"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.24"
$ws.Range("E40").Value = "  -5.58%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.44"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "2.470.18"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.98"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.70"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -1.21%  "
